# CIMS_Deliverables.pptx edit: mark "Connect Database" and "Create Login
# Page" tasks as COMPLETE on the Project Plan table (slide with
# SlideID=290, the "Project Plan" table slide).
#
# Commit message: "added login screen & database connection"
# -> the login screen and database connection work items are now done,
#    so their Status cells get the "COMPLETE" label (matching the
#    formatting already used for the other completed "Build Database"
#    row: 19pt text).

$p = $ppt.ActivePresentation

# Locate the slide by its persistent SlideID (290) rather than trusting a
# fixed index, in case slide order ever shifts.
$targetSlide = $null
foreach ($sl in $p.Slides) {
    if ($sl.SlideID -eq 290) {
        $targetSlide = $sl
        break
    }
}

$tbl = $targetSlide.Shapes.Item("Table 1").Table

# Column 4 is "Status". Row 4 = "Connect Database" (Phase 1), row 5 =
# "Create Login Page" (Phase 2) - both currently have an empty Status
# cell and need to read "COMPLETE", sized like the existing completed
# entry (19pt).
$statusCol = 4

$connectDatabaseRow = 4
$createLoginPageRow = 5

$cell = $tbl.Cell($connectDatabaseRow, $statusCol)
$cell.Shape.TextFrame.TextRange.Text = "COMPLETE"
$cell.Shape.TextFrame.TextRange.Font.Size = 19

$cell = $tbl.Cell($createLoginPageRow, $statusCol)
$cell.Shape.TextFrame.TextRange.Text = "COMPLETE"
$cell.Shape.TextFrame.TextRange.Font.Size = 19
